# Apply translation-sheet reorganisation:
#  - remove the two duplicate 'ACORN data is not of the right format...' rows,
#    the old 'HAI point prevalence by type of ward' row and the old
#    'Select lab data format:' row from their original (sorted) position and
#    append them at the end of the table marked status='deleted'
#  - insert four new source strings (status='new', translation='TBT'):
#    'Contains names of organisms before and after mapping.',
#    'Download Lab Log (.xlsx)', 'HAI point prevalence by ' and
#    "Remove 'Not Cultured' specimens"
#  - this grows the sheet from 181 to 185 data/header rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full A2:C185 block (en / la / status) in its final row order
$data = New-Object 'object[,]' 184,3
$data[0,0] = '.acorn data successfully generated!'
$data[0,1] = 'ສຳເລັດການສ້າງຂໍ້ມູນ.acorn!'
$data[0,2] = ''
$data[1,0] = '.acorn file saved on server.'
$data[1,1] = 'ເອກະສານ.acorn ບັນທຶກລົງໃນ server ແລ້ວ'
$data[1,2] = ''
$data[2,0] = '.acorn not saved.'
$data[2,1] = '.acorn ບໍ່ໄດ້ຖືກບັນທຶກ.'
$data[2,2] = ''
$data[3,0] = '(1/4) Download Clinical data'
$data[3,1] = '(1/4) ດາວໂລດຂໍ້ມູນທາງຄຣີນິກ'
$data[3,2] = ''
$data[4,0] = '(2/4) Provide Lab data'
$data[4,1] = '(2/4) ຈັດຫາຂໍ້ມູນແລັບ'
$data[4,2] = ''
$data[5,0] = '(3/4) Combine Clinical and Lab data'
$data[5,1] = '(3/4) ລວມຂໍ້ມູນຄຣີນິກ ແລະ ຂໍ້ມູນແລັບເຂົ້າກັນ'
$data[5,2] = ''
$data[6,0] = '(4/4) Save .acorn file'
$data[6,1] = '(4/4)ບັນທຶກເອກະສານຂອງ.acorn'
$data[6,2] = ''
$data[7,0] = '(Optional) Comments:'
$data[7,1] = '(ຕົວເລືອກ) ຄຳຄິດເຫັນ'
$data[7,2] = ''
$data[8,0] = '(To log out, close the app.)'
$data[8,1] = 'TBT'
$data[8,2] = ''
$data[9,0] = 'ACORN Participating Countries'
$data[9,1] = 'ປະເທດທີ່ເຂົ້າຮ່ວມ ACORN'
$data[9,2] = ''
$data[10,0] = 'All ''orgname'' are provided.'
$data[10,1] = 'ທຸກໆ ''orgname''ທັງໝົດແມ່ນຖືກຕອບສະໜອງ'
$data[10,2] = ''
$data[11,0] = 'All ''patid'' are provided.'
$data[11,1] = 'patid'' ທັງໝົດແມ່ນຖືກຕອບສະໜອງ'
$data[11,2] = ''
$data[12,0] = 'All ''specdate'' are provided.'
$data[12,1] = 'specdate'' ທັງໝົດແມ່ນຖືກຕອບສະໜອງ'
$data[12,2] = ''
$data[13,0] = 'All ''specdate'' are today or before today.'
$data[13,1] = 'ທຸກໆ ''specdate'' ແມ່ນມື້ນີ້ ຫຼື ມື້ກ່ອນໜ້ານີ້'
$data[13,2] = ''
$data[14,0] = 'All ''specgroup'' are provided.'
$data[14,1] = 'ທັງໝົດຂອງ''specgroup'' ແມ່ນຖືກກຳນົດໄວ້ແລ້ວ'
$data[14,2] = ''
$data[15,0] = 'All ''specid'' are provided.'
$data[15,1] = 'specid'' ທັງໝົດແມ່ນຖືກຕອບສະໜອງ'
$data[15,2] = ''
$data[16,0] = 'All dates of enrolment for HAI patients have a matching date in the HAI survey dataset'
$data[16,1] = 'ທຸກວັນທີເຂົ້າການສຶກສາສຳລັບຄົນເຈັບHAI ຈັບຄູ່ເຂົ້າກັບວັນທີເຮັດການສຳຫຼວດ HAI'
$data[16,2] = ''
$data[17,0] = 'All Other Organisms'
$data[17,1] = 'ຕົວເຊື້ອອື່ນໆທັງໝົດ'
$data[17,2] = ''
$data[18,0] = 'All valid records have an ACORN ID.'
$data[18,1] = 'TBT'
$data[18,2] = ''
$data[19,0] = 'AMR'
$data[19,1] = 'TBT'
$data[19,2] = ''
$data[20,0] = 'and generate enrolment log.'
$data[20,1] = 'ແລະສ້າງບັນທຶກການເຂົ້າຮ່ວມ'
$data[20,2] = ''
$data[21,0] = 'Attempting to connect.'
$data[21,1] = 'ກຳລັງພະຍາຍາມເຊື່ອມຕໍ່'
$data[21,2] = ''
$data[22,0] = 'Blood culture collected within 24 hours of admission (CAI) / symptom onset (HAI)'
$data[22,1] = 'ປູກເລືອດພາຍໃນ24ຊົ່ວໂມງຂອງການເຂົ້ານອນ (CAI) / ເລີ້ມມີການຊຶມເຊື້ອໃນໂຮງໝໍ (HAI)'
$data[22,2] = ''
$data[23,0] = 'Blood Culture Contaminants'
$data[23,1] = 'ການປົນເປື້ອນຂອງການປູກເລືອດ'
$data[23,2] = ''
$data[24,0] = 'Bloodstream Infection (BSI)'
$data[24,1] = 'ການຊຶມເຊື້ອໃນກະແສເລືອດ'
$data[24,2] = ''
$data[25,0] = 'Calculated age is consistent with ''Age Category'''
$data[25,1] = 'ການຄິດໄລ່ອາຍຸແມ່ນເໝາະສົມກັບ ''Age Category'''
$data[25,2] = ''
$data[26,0] = 'Calculated age isn''t always consistent with ''Age Category'''
$data[26,1] = 'ການຄິດໄລ່ອາຍຸອາດບໍ່ເໝາະສົມກັບ ''Age Category'' ສະເໝີໄປ'
$data[26,2] = ''
$data[27,0] = 'Cancel'
$data[27,1] = 'ຍົກເລີກ'
$data[27,2] = ''
$data[28,0] = 'Care should be taken when interpreting rates and AMR profiles where there are small numbers of cases or bacterial isolates: point estimates may be unreliable.'
$data[28,1] = 'ຄວນໃຊ້ຄວາມລະມັດລະວັງໃນການແປຄວາມໝາຍອັດຕາແລະຮູບຮ່າງລັກສະນະຂອງການຕ້ານຕໍ່ຢາຕ້ານເຊື້ອທີ່ມີຈຳນວນຄົນເຈັບໜ້ອຍ ຫຼື ແບັກທີເຣຍທີ່ແຍກໄດ້: ການແປຜົນອາດບໍ່ໜ້າເຊື່ອຖື'
$data[28,2] = ''
$data[29,0] = 'Clinical and day-28 outcomes are consistent.'
$data[29,1] = 'ອາການຄຣີນິກ ແລະ ການຕິດຕາມມື້ທີ 28 ແມ່ນກົງກັນ'
$data[29,2] = ''
$data[30,0] = 'Clinical and day-28 outcomes aren''t consistent for some dead patients.'
$data[30,1] = 'ອາການຄຣີນິກ ແລະ ການຕິດຕາມມື້ທີ 28 ບໍ່ກົງກັນໃນບາງຄົນເຈັບທີ່ເສຍຊີວິດ'
$data[30,2] = ''
$data[31,0] = 'Clinical Outcome'
$data[31,1] = 'ອາການຫຼັງຕິດຕາມ'
$data[31,2] = ''
$data[32,0] = 'Clinical Outcome Status:'
$data[32,1] = 'ສະຖານະພາບອາການຄຣີນິກຫຼັງຕິດຕາມ'
$data[32,2] = ''
$data[33,0] = 'Co-resistances'
$data[33,1] = 'TBT'
$data[33,2] = ''
$data[34,0] = 'Combine Susceptible + Intermediate'
$data[34,1] = 'Susceptible + Intermediate ລວມເຂົ້າກັນ'
$data[34,2] = ''
$data[35,0] = 'Consider saving .acorn file on the cloud for additional security.'
$data[35,1] = 'ຄວນພິຈາລະນາບັນທຶກຂໍ້ມູນ.acorn ໃນ cloud ເພື່ອຄວາມປອດໄພເພີ້ມເຕີມ'
$data[35,2] = ''
$data[36,0] = 'Contains names of organisms before and after mapping.'
$data[36,1] = 'TBT'
$data[36,2] = 'new'
$data[37,0] = 'Couldn''t connect to server. Please check internet access.'
$data[37,1] = 'ບໍ່ສາມາດເຊື່ອມຕໍ່ກັບ server ໄດ້. ກະລຸນາກວດສອບການເຂົ້າເຖິງອີນເຕີເນັດ'
$data[37,2] = ''
$data[38,0] = 'Critical errors with clinical data.'
$data[38,1] = 'ຂໍ້ຜິດພາດທີ່ຮ້າຍແຮງກ່ຽວກັບຂໍ້ມູນທາງຄຣີນິກ.'
$data[38,2] = ''
$data[39,0] = 'Culture results per specimen type'
$data[39,1] = 'ຜົນການປູກຕໍ່ກັບປະເພດຕົວຢ່າງ'
$data[39,2] = ''
$data[40,0] = 'Data Management'
$data[40,1] = 'ການຈັດການຂໍ້ມູນ'
$data[40,2] = ''
$data[41,0] = 'Date of Enrolment'
$data[41,1] = 'ວັນທີເຂົ້າຮ່ວມການສຶກສາ'
$data[41,2] = ''
$data[42,0] = 'Day 28'
$data[42,1] = 'ມື້ທີ28'
$data[42,2] = ''
$data[43,0] = 'Day 28 Status:'
$data[43,1] = 'ສະຖານະພາບມື້ທີ28'
$data[43,2] = ''
$data[44,0] = 'Diagnosis at Enrolment'
$data[44,1] = 'ການບົ່ງມະຕິເວລາເຂົ້າການສຶກສາ'
$data[44,2] = ''
$data[45,0] = 'Dismiss'
$data[45,1] = 'TBT'
$data[45,2] = ''
$data[46,0] = 'Distribution of Enrolments'
$data[46,1] = 'ການແຈກຢາຍຂອງການເຂົ້າຮ່ວມ'
$data[46,2] = ''
$data[47,0] = 'Download Enrolment Log (.xlsx)'
$data[47,1] = 'ດຶງຂໍ້ມູນບັນທຶກການເຂົ້າຮ່ວມ (.xlsx)'
$data[47,2] = ''
$data[48,0] = 'Download Lab Log (.xlsx)'
$data[48,1] = 'TBT'
$data[48,2] = 'new'
$data[49,0] = 'Empiric Antibiotics Prescribed'
$data[49,1] = 'ການໃຫ້ຢາຕ້ານເຊື້ອກ່ອນການບົ່ງມະຕິຊັດເຈນ'
$data[49,2] = ''
$data[50,0] = 'Enrolments'
$data[50,1] = 'ເຂົ້າການສືກສາ'
$data[50,2] = ''
$data[51,0] = 'Enrolments by (type of) Ward'
$data[51,1] = 'ເຂົ້າຮ່ວມການສຶກສາໂດຍ(ປະເພດຂອງ)ພະແນກ'
$data[51,2] = ''
$data[52,0] = 'Enrolments with Blood Culture'
$data[52,1] = 'ເຂົ້າຮ່ວມການສຶກສາສົມທົບກັບປູກເລືອດ'
$data[52,2] = ''
$data[53,0] = 'Error in combining clinical and lab data.'
$data[53,1] = 'ເກີດຂໍ້ຜິດພາດໃນການລວມຂໍ້ມູນທາງຄຣີນິກກັບຂໍ້ມູນທາງຫ້ອງວິເຄາະ'
$data[53,2] = ''
$data[54,0] = 'Every D28 record (F04) matches exactly one patient enrolment (F01).'
$data[54,1] = 'ທຸກໆຟອມບັນທຶກມື້ທີ28 (F04) ຈັບຄູ່ກັບຟອມເຂົ້າການສືກສາ (F01).'
$data[54,2] = ''
$data[55,0] = 'Every hospital outcome record (F03) has a matching infection episode (F02).'
$data[55,1] = 'ທຸກໆຟອມບັນທຶກອອກໂຮງໝໍ (F03) ຈັບຄູ່ເຂົ້າກັບຟອມພາກການຊຶມເຊື້ອ (F02).'
$data[55,2] = ''
$data[56,0] = 'Every hospital outcome record (F03) has a matching patient enrolment (F01).'
$data[56,1] = 'ບາງຟອມບັນທຶກອອກໂຮງໝໍ (F03) ຈັບຄູ່ເຂົ້າກັບຟອມເຂົ້າຮ່ວມການສຶກສາ (F01)'
$data[56,2] = ''
$data[57,0] = 'Every infection episode record (F02) has a matching patient enrolment (F01).'
$data[57,1] = 'ທຸກໆຟອມບັນທຶກພາກການຊຶມເຊື້ອ (F02) ຈັບຄູ່ເຂົ້າກັບຟອມເຂົ້າການສຶກສາ (F01)'
$data[57,2] = ''
$data[58,0] = 'File name:'
$data[58,1] = 'ຊື່ເອກະສານ'
$data[58,2] = ''
$data[59,0] = 'First sheet is the log of all enrolments retrived from REDCap (as per adjacent table). The second sheet is a listing of all flagged elements.'
$data[59,1] = 'ໜ້າທຳອິດເປັນບັນທຶກການລົງທະບຽນເຂົ້າຮ່ວມທັງໝົດທີ່ດຶງມາຈາກREDCap (ຕາມຕາຕະລາງທີ່ຢູ່ໃກ້ກັນ) ໜ້າທີສອງແມ່ນການສະແດງລາຍການອົງປະກອບທີ່ຖືກຕັ້ງຄ່າແລ້ວ'
$data[59,2] = ''
$data[60,0] = 'Follow-up'
$data[60,1] = 'ການຕິດຕາມ'
$data[60,2] = ''
$data[61,0] = 'from cultures that have growth'
$data[61,1] = 'ມີເຊື້ອເກີດຈາກການປູກ'
$data[61,2] = ''
$data[62,0] = 'Generate .acorn file'
$data[62,1] = 'ສ້າງເອກະສານຂອງ.acorn'
$data[62,2] = ''
$data[63,0] = 'Generate and load .acorn from clinical and lab data'
$data[63,1] = 'TBT'
$data[63,2] = ''
$data[64,0] = 'Generating .acorn'
$data[64,1] = 'ກຳລັງສ້າງ .acorn'
$data[64,2] = ''
$data[65,0] = 'Get data from REDCap'
$data[65,1] = 'ເອົາຂໍ້ມູນຈາກ REDCap'
$data[65,2] = ''
$data[66,0] = 'Get the latest production release'
$data[66,1] = 'TBT'
$data[66,2] = ''
$data[67,0] = 'Growth / No Growth'
$data[67,1] = 'ເກີດເຊື້ອ/ບໍ່ເກີດເຊື້ອ'
$data[67,2] = ''
$data[68,0] = 'HAI point prevalence by '
$data[68,1] = 'TBT'
$data[68,2] = 'new'
$data[69,0] = 'HAI Prevalence'
$data[69,1] = 'ຄວາມຊຸກຊຸມຂອງການຊຶມເຊື້ອໃນໂຮງໝໍ'
$data[69,2] = ''
$data[70,0] = 'Horizontal bars show the size of a set of SR results while vertical bars show the number of resistant isolates for the corresponding antibiotic.'
$data[70,1] = 'TBT'
$data[70,2] = ''
$data[71,0] = 'Info on loaded .acorn'
$data[71,1] = 'TBT'
$data[71,2] = ''
$data[72,0] = 'Initial & Final Surveillance Diagnosis'
$data[72,1] = 'ການບົ່ງມະຕິການເຝົ້າລະວັງໃນເບື້ອງຕົ້ນແລະຂັ້ນສຸດທ້າຍ'
$data[72,2] = ''
$data[73,0] = 'Isolates'
$data[73,1] = 'ການແຍກເຊື້ອ'
$data[73,2] = ''
$data[74,0] = 'Issue detected with REDCap data. Please report to ACORN data managers. Until resolution, only existing .acorn files can be used.'
$data[74,1] = 'TBT'
$data[74,2] = ''
$data[75,0] = 'It might take a couple of minutes. This window will close on completion.'
$data[75,1] = 'ອາດໃຊ້ເວລາເລັກນ້ອຍ, ໜ້າຕ່າງນີ້ຈະປິດລົງເມື່ອສຳເລັດ'
$data[75,2] = ''
$data[76,0] = 'Lab data successfully processed!'
$data[76,1] = 'ສຳເລັດການປະມວນຜົນຂໍ້ມູນທາງຫ້ອງວິເຄາະ!'
$data[76,2] = ''
$data[77,0] = 'Lab data successfully provided.'
$data[77,1] = 'ໃຫ້ຂໍ້ມູນແລັບສຳເລັດແລ້ວ'
$data[77,2] = ''
$data[78,0] = 'Lab dataset contains the minimal columns.'
$data[78,1] = 'ໃນຊຸດຂໍ້ມູນຂອງແລັບມີຖັນຈຳກັດຈຳນວນໜ້ອຍ'
$data[78,2] = ''
$data[79,0] = 'Lab dataset does not contains the minimal columns.'
$data[79,1] = 'ໃນຊຸດຂໍ້ມູນຂອງແລັບບໍ່ມີຖັນຈຳກັດຈຳນວນໜ້ອຍ'
$data[79,2] = ''
$data[80,0] = 'Language'
$data[80,1] = 'ພາສາ'
$data[80,2] = ''
$data[81,0] = 'Load .acorn'
$data[81,1] = 'Load .acorn'
$data[81,2] = ''
$data[82,0] = 'Load .acorn from cloud'
$data[82,1] = 'TBT'
$data[82,2] = ''
$data[83,0] = 'Load .acorn from local file'
$data[83,1] = 'TBT'
$data[83,2] = ''
$data[84,0] = 'Load selected .acorn'
$data[84,1] = 'ເລືອກ Load .acorn'
$data[84,2] = ''
$data[85,0] = 'Loading data.'
$data[85,1] = 'ກຳລັງດຶງຂໍ້ມູນ'
$data[85,2] = ''
$data[86,0] = 'Log in'
$data[86,1] = 'ເຂົ້າສູ່ລະບົບ'
$data[86,2] = ''
$data[87,0] = 'Microbiology'
$data[87,1] = 'ຈຸລີນຊີວິທະຍາ'
$data[87,2] = ''
$data[88,0] = 'Most frequent 10 organisms in the plot and complete listing in the table. Contaminants are in red.'
$data[88,1] = '10 ຕົວເຊື້ອທີ່ມັກພົບຫຼາຍທີ່ສຸດຖືກສະແດງໃນຮູບພາບ ແລະ ລາຍການທັງໝົດແມ່ນຖືກນຳສະເໜີຢູ່ໃນຕາຕະລາງ. ເຊື້ອປົນເປື້ອນແມ່ນເປັນສີແດງ'
$data[88,2] = ''
$data[89,0] = 'No .acorn data loaded.'
$data[89,1] = 'ບໍ່ມີຂໍ້ມູນ.acorn ທີ່ຖືກດຶງ.'
$data[89,2] = ''
$data[90,0] = 'No Blood Culture'
$data[90,1] = 'ບໍ່ມີການເລືອດປູກ'
$data[90,2] = ''
$data[91,0] = 'Not connected to internet.'
$data[91,1] = 'ບໍ່ໄດ້ເຊື່ອມຕໍ່ອີນເຕີເນັດ.'
$data[91,2] = ''
$data[92,0] = 'Number of specimens per specimen type'
$data[92,1] = 'ຈຳນວນຕົວຢ່າງຕໍ່ກັບປະເພດຂອງຕົວຢ່າງ'
$data[92,2] = ''
$data[93,0] = 'Occupancy rate per type of ward per month'
$data[93,1] = 'ອັດຕາການເຂົ້ານອນຕໍ່ປະເພດຂອງພະແນກຕໍ່ເດືອນ'
$data[93,2] = ''
$data[94,0] = 'of blood cultures grew a potential contaminant.'
$data[94,1] = 'ຂອງການປູກເລືອດເປັນໄປໄດ້ອາດເກີດເຊື້ອປົນເປື້ອນ'
$data[94,2] = ''
$data[95,0] = 'of cultures have growth.'
$data[95,1] = 'ຂອງການປູກມີເກີດເຊື້ອ'
$data[95,2] = ''
$data[96,0] = 'of enrolments with blood culture.'
$data[96,1] = 'ຂອງການເຂົ້າການສຶກສາກັບການປູກເລືອດ'
$data[96,2] = ''
$data[97,0] = 'of Target Pathogens'
$data[97,1] = 'ຂອງເຊື້ອເປົ້າໝາຍ'
$data[97,2] = ''
$data[98,0] = 'Only isolates that have been tested against all of the drugs are included in the upset plot.'
$data[98,1] = 'TBT'
$data[98,2] = ''
$data[99,0] = 'Overview'
$data[99,1] = 'ພາບລວມ'
$data[99,2] = ''
$data[100,0] = 'Password'
$data[100,1] = 'ລະຫັດຜ່ານ'
$data[100,2] = ''
$data[101,0] = 'Patient Age Distribution'
$data[101,1] = 'ການແຈກຢາຍອາຍຸຂອງຄົນເຈັບ'
$data[101,2] = ''
$data[102,0] = 'Patient Comorbidities'
$data[102,1] = 'ພະຍາດປະຈຳໂຕຂອງຄົນເຈັບ'
$data[102,2] = ''
$data[103,0] = 'Patient enrolments'
$data[103,1] = 'ຄົນເຈັບທີ່ເຂົ້າການສຶກສາ'
$data[103,2] = ''
$data[104,0] = 'Patients Transferred'
$data[104,1] = 'ຄົນເຈັບຖືກນຳສົ່ງ'
$data[104,2] = ''
$data[105,0] = 'Please log in'
$data[105,1] = 'ກະລຸນາລົງທະບຽນເຂົ້າສູ່ລະບົບ'
$data[105,2] = ''
$data[106,0] = 'Processing lab data.'
$data[106,1] = 'ກຳລັງປະມວນຜົນຂໍ້ມູນຫ້ອງວິເຄາະ'
$data[106,2] = ''
$data[107,0] = 'Reading lab data.'
$data[107,1] = 'ກຳລັງອ່ານຂໍ້ມູນຫ້ອງວິເຄາະ'
$data[107,2] = ''
$data[108,0] = 'Remove ''Not Cultured'' specimens'
$data[108,1] = 'TBT'
$data[108,2] = 'new'
$data[109,0] = 'Remove blood culture contaminants from the following visualizations'
$data[109,1] = 'ການປູກເລືອດທີ່ມີການປົນເປື້ອນແມ່ນລົບຜົນອອກຈາກຜົນການສະແດງຂໍ້ມູນ'
$data[109,2] = ''
$data[110,0] = 'Reset Enrolments Filters'
$data[110,1] = 'Reset Enrolments Filters'
$data[110,2] = ''
$data[111,0] = 'Resistance to 3rd gen. Cephalosporins Over Time'
$data[111,1] = 'ການຕ້ານຕໍ່ 3rd gen. Cephalosporins ຄ່ອຍເປັນຄ່ອຍໄປຕາມໄລຍະເວລາ'
$data[111,2] = ''
$data[112,0] = 'Resistance to Carbapenems Over Time'
$data[112,1] = 'ມີການຕ້ານຕໍ່ຢາ Carbapenems ຄ່ອຍເປັນຄ່ອຍໄປຕາມໄລຍະເວລາ'
$data[112,2] = ''
$data[113,0] = 'Resistance to Fluoroquinolones Over Time'
$data[113,1] = 'ການຕ້ານຕໍ່ Fluoroquinolones ຄ່ອຍເປັນຄ່ອຍໄປຕາມໄລຍະເວລາ'
$data[113,2] = ''
$data[114,0] = 'Resistance to Oxacillin Over Time'
$data[114,1] = 'ການຕ້ານຕໍ່ຢາ Oxacillin ເປັນໄປຕາມໄລຍະເວລາ'
$data[114,2] = ''
$data[115,0] = 'Resistance to Penicillin G - meningitis Over Time'
$data[115,1] = 'ການຕ້ານຕໍ່ຢາ Penicillin G- meningitis ເປັນໄປຕາມໄລຍະເວລາ'
$data[115,2] = ''
$data[116,0] = 'Resistance to Penicillin G Over Time'
$data[116,1] = 'ການຕ້ານຕໍ່ຢາ Penicillin G ເປັນໄປຕາມໄລຍະເວລາ'
$data[116,2] = ''
$data[117,0] = 'Retriving data from REDCap server.'
$data[117,1] = 'ກຳລັງກູ້ຄືນຂໍ້ມູນຈາກຖານ REDCap'
$data[117,2] = ''
$data[118,0] = 'Save .acorn file'
$data[118,1] = 'ບັນທຶກເອກະສານຂອງ.acorn'
$data[118,2] = ''
$data[119,0] = 'Save acorn data'
$data[119,1] = 'ບັນທຶກຂໍ້ມູນ acorn'
$data[119,2] = ''
$data[120,0] = 'Save on Server'
$data[120,1] = 'ບັນທຶກລົງໃນຖານຂໍ້ມູນ'
$data[120,2] = ''
$data[121,0] = 'See Breakdown by Ward'
$data[121,1] = 'ເບິ່ງລາຍລະອຽດໂດຍອິງໃສ່ພະແນກ'
$data[121,2] = ''
$data[122,0] = 'See by Week'
$data[122,1] = 'ເບິງລາຍອາທິດ'
$data[122,2] = ''
$data[123,0] = 'Show antibiotics combinations'
$data[123,1] = 'TBT'
$data[123,2] = ''
$data[124,0] = 'Show comorbidities combinations'
$data[124,1] = 'ສະແດງພະຍາດປະຈຳໂຕອື່ນຮ່ວມ'
$data[124,2] = ''
$data[125,0] = 'SIR Evaluation'
$data[125,1] = 'TBT'
$data[125,2] = ''
$data[126,0] = 'Some D28 records (F04) don''t have a matching patient enrolment (F01).'
$data[126,1] = 'ບາງຟອມບັນທຶກມື້ທີ28 (F04 ) ບໍ່ຈັບຄູ່ກັບຟອມເຂົ້າການສຶກສາ (F01).'
$data[126,2] = ''
$data[127,0] = 'Some dates of enrolment for HAI patients do have a matching date in the HAI survey dataset'
$data[127,1] = 'TBT'
$data[127,2] = ''
$data[128,0] = 'Some hospital outcome records (F03) don''t have a matching infection episode (F02). These records have been removed.'
$data[128,1] = 'ບາງຟອມບັນທຶກອອກໂຮງໝໍ (F03) ບໍ່ຈັບຄູ່ເຂົ້າກັບຟອມພາກການຊຶມເຊື້ອ (F02). ຟອມເຫຼົ່ານີ້ຖືກລຶບແລ້ວ'
$data[128,2] = ''
$data[129,0] = 'Some hospital outcome records (F03) don''t have a matching patient enrolment (F01).'
$data[129,1] = 'ບາງຟອມບັນທຶກອອກໂຮງໝໍ (F03) ບໍ່ຈັບຄູ່ເຂົ້າກັບແບບຟອມເຂົ້າການສຶກສາ (F01)'
$data[129,2] = ''
$data[130,0] = 'Some infection episode records (F02) don''t have a matching patient enrolment (F01). These records have been removed.'
$data[130,1] = 'ບາງຟອມບັນທຶກພາກການຊຶມເຊື້ອ (F02) ບໍ່ຈັບຄູ່ເຂົ້າກັບຟອມເຂົ້າການສຶກສາ (F01).ບັນທຶກເຫຼົ່ານີ້ຖືກລຶບແລ້ວ'
$data[130,2] = ''
$data[131,0] = 'Some records with a missing ACORN ID. These records have been removed.'
$data[131,1] = 'TBT'
$data[131,2] = ''
$data[132,0] = 'Specimen Types'
$data[132,1] = 'ປະເພດຂອງຕົວຢ່າງ'
$data[132,2] = ''
$data[133,0] = 'Specimens'
$data[133,1] = 'TBT'
$data[133,2] = ''
$data[134,0] = 'Specimens Collected'
$data[134,1] = 'ການເກັບຕົວຢ່າງ'
$data[134,2] = ''
$data[135,0] = 'specimens per enrolment'
$data[135,1] = 'ຕົວຢ່າງຕໍ່ການເຂົ້າຮ່ວມການສຶກສາ'
$data[135,2] = ''
$data[136,0] = 'Successfully combined clinical and lab data into .acorn file'
$data[136,1] = 'ສຳເລັດການລວມຂໍ້ມູນທາງຄຣີນິກກັບຂໍ້ມູນທາງຫ້ອງວິເຄາະເຂົ້າໃນ .acorn file'
$data[136,2] = ''
$data[137,0] = 'Successfully loaded data.'
$data[137,1] = 'ສຳເລັດການດຶງຂໍ້ມູນແລ້ວ'
$data[137,2] = ''
$data[138,0] = 'Successfully logged in.'
$data[138,1] = 'ສຳເລັດການລົງທະບຽນເຂົ້າສູ່ລະບົບ'
$data[138,2] = ''
$data[139,0] = 'Successfully saved .acorn file in the cloud. You can now explore acorn data.'
$data[139,1] = 'ສຳເລັດການບັນທຶກເອກະສານ .acorn ໃນ cloud. ທ່ານສາມາດກວດສອບຂໍ້ມູນ .acorn ໄດ້ແລ້ວ'
$data[139,2] = ''
$data[140,0] = 'Successfully saved .acorn file locally.'
$data[140,1] = 'ສຳເລັດການບັນທຶກເອກກະສານ.acorn ໃນເຄື່ອງ'
$data[140,2] = ''
$data[141,0] = 'Supply first valid clinical and lab data.'
$data[141,1] = 'ຈັດຫາຂໍ້ມູນຄຣີນິກແລະຫ້ອງວິເຄາະຄົບຖ້ວນເປັນອັນດັບທຳອິດ'
$data[141,2] = ''
$data[142,0] = 'Susceptible & Intermediate are always combined in this visualisation of co-resistances.'
$data[142,1] = 'TBT'
$data[142,2] = ''
$data[143,0] = 'The 10 most common initial-final diagnosis combinations:'
$data[143,1] = '10 ອັນດັບການບົ່ງມະຕິໃນເບື້ອງຕົ້ນ ແລະ ສຸດທ້າຍ ທີ່ພົບຫຼາຍທີ່ສຸດ'
$data[143,2] = ''
$data[144,0] = 'The following ''patient id'' are atypical cases (one HCAI/CAI with early HAI but no overlap):'
$data[144,1] = 'TBT'
$data[144,2] = ''
$data[145,0] = 'The following ''patient id'' are problem case (overlapping specimen collection windows):'
$data[145,1] = 'ການເຊື່ອມຕໍ່ລະຫັດຂໍ້ມູນຄົນເຈັບເກີດຂໍ້ຜິດພາດ (ເກີດການທັບຊ້ອນຂໍ້ມູນຕົວຢ່າງທີ່ເກັບລວບລວມ):'
$data[145,2] = ''
$data[146,0] = 'The REDCap dataset is empty/in wrong format. Please contact ACORN support.'
$data[146,1] = 'TBT'
$data[146,2] = ''
$data[147,0] = 'The REDCap dataset is in the right format.'
$data[147,1] = 'TBT'
$data[147,2] = ''
$data[148,0] = 'There are D28 follow-up done before the expected D28 date.'
$data[148,1] = 'TBT'
$data[148,2] = ''
$data[149,0] = 'There are multiple F02 with identical ACORN ID, admission date, and episode enrolment date.'
$data[149,1] = 'TBT'
$data[149,2] = ''
$data[150,0] = 'There are no atypical case (one HCAI/CAI with early HAI but no overlap).'
$data[150,1] = 'TBT'
$data[150,2] = ''
$data[151,0] = 'There are no D28 follow-up done before the expected D28 date.'
$data[151,1] = 'TBT'
$data[151,2] = ''
$data[152,0] = 'There are no isolate with valid AST results. Please contact ACORN support.'
$data[152,1] = 'ບໍ່ມີຜົນແຍກການທົດສອບຢາຕ້ານເຊື້ອທີ່ສົມບູນ. ກະລຸນາຕິດຕໍ່ຝ່າຍສະໜັບສະໜຸນ ACORN'
$data[152,2] = ''
$data[153,0] = 'There are no multiple F02 with identical ACORN ID, admission date, and episode enrolment date.'
$data[153,1] = 'TBT'
$data[153,2] = ''
$data[154,0] = 'There are no problem case (overlapping specimen collection windows)'
$data[154,1] = 'ບໍ່ພົບກໍລະນີທີ່ມີບັນຫາ ( ການທັບຊ້ອນຂໍ້ມູນຕົວຢ່າງທີ່ເກັບລວບລວມ):'
$data[154,2] = ''
$data[155,0] = 'There are rows for which ''specdate'' are after today.'
$data[155,1] = 'ມີແຖວຂອງ''specdate'' ເເມ່ນມື້ຫຼັງຈາກນີ້'
$data[155,2] = ''
$data[156,0] = 'There are rows with missing ''orgname''.'
$data[156,1] = 'ບາງແຖວບໍ່ມີ ''orgname''.'
$data[156,2] = ''
$data[157,0] = 'There are rows with missing ''patid''.'
$data[157,1] = 'ບາງແຖວບໍ່ມີ ''patid'''
$data[157,2] = ''
$data[158,0] = 'There are rows with missing ''specdate''.'
$data[158,1] = 'ບາງແຖວບໍ່ມີ ''specdate'''
$data[158,2] = ''
$data[159,0] = 'There are rows with missing ''specgroup''.'
$data[159,1] = 'ບາງແຖວບໍ່ມີ ''specgroup''.'
$data[159,2] = ''
$data[160,0] = 'There are rows with missing ''specid''.'
$data[160,1] = 'ບາງແຖວບໍ່ມີ ''specid'''
$data[160,2] = ''
$data[161,0] = 'There is a critical issue with clinical data. The issue should be fixed in REDCap.'
$data[161,1] = 'ມີບັນຫາຮ້າຍແຮງກ່ຽວກັບຂໍ້ມູນທາງຄຣີນິກ. ຄວນໄດ້ຮັບການແກ້ໄຂຢູ່ໃນ REDCap.'
$data[161,2] = ''
$data[162,0] = 'There is no data to display for this organism.'
$data[162,1] = 'ບໍ່ມີຂໍ້ມູນນຳສະເໜີສຳລັບຕົວເຊື້ອນີ້'
$data[162,2] = ''
$data[163,0] = 'There is no HAI survey data'
$data[163,1] = 'ບໍ່ມີຂໍ້ມູນການສຳຫຼວດ HAI'
$data[163,2] = ''
$data[164,0] = 'Trying to save .acorn file on server.'
$data[164,1] = 'ກຳລັງພະຍາຍຍາມບັນທຶກເອກະສານ .acorn ລົງ server'
$data[164,2] = ''
$data[165,0] = 'Updated Charlson Comorbidity Index (uCCI)'
$data[165,1] = 'TBT'
$data[165,2] = ''
$data[166,0] = 'User'
$data[166,1] = 'ຜູ້ໃຊ້ງານ'
$data[166,2] = ''
$data[167,0] = 'Variables in Table:'
$data[167,1] = 'ຕົວແປໃນຕາຕະລາງ'
$data[167,2] = ''
$data[168,0] = 'Ward Occupancy Rates'
$data[168,1] = 'ອັດຕາຄວາມໜາແໜ້ນໃນພະແນກ'
$data[168,2] = ''
$data[169,0] = 'We couldn''t download the lab codes file. Please contact ACORN support.'
$data[169,1] = 'ພວກເຮົາບໍ່ສາມາດດາວໂລດເອກະສານລະຫັດແລັບ. ກະລຸນາຕິດຕໍ່ຜູ້ສະໜັບສະໜຸນ ACORN'
$data[169,2] = ''
$data[170,0] = 'We couldn''t download the lab data dictionary. Please contact ACORN support'
$data[170,1] = 'TBT'
$data[170,2] = ''
$data[171,0] = 'Welcome'
$data[171,1] = 'ຍີນດີຕ້ອນຮັບ'
$data[171,2] = ''
$data[172,0] = 'What do you want to do?'
$data[172,1] = 'ທ່ານຕ້ອງການຢາກເຮັດຫຍັງ?'
$data[172,2] = ''
$data[173,0] = 'With Microbiology'
$data[173,1] = 'ຮ່ວມກັບຈຸລີນຊີວິທະຍາ'
$data[173,2] = ''
$data[174,0] = 'Wrong connection credentials.'
$data[174,1] = 'ຂໍ້ມູນຮັບຮອງການເຊື່ອມຕໍ່ບໍ່ຖືກຕ້ອງ'
$data[174,2] = ''
$data[175,0] = 'You are running ACORN dashboard'
$data[175,1] = 'TBT'
$data[175,2] = ''
$data[176,0] = 'You can check here if it''s the latest production release.'
$data[176,1] = 'TBT'
$data[176,2] = ''
$data[177,0] = 'Your ACORN dashboard is up to date'
$data[177,1] = 'TBT'
$data[177,2] = ''
$data[178,0] = 'Follow us on Twitter'
$data[178,1] = 'ຕິດຕາມພວກເຮົາໄດ້ທີ່ Twitter'
$data[178,2] = ''
$data[179,0] = 'Records in Lab data and BSI forms:'
$data[179,1] = 'ບັນທຶກຂໍ້ມູນໃນຫ້ອງວິເຄາະ ແລະ ແບບຟອມ BSI'
$data[179,2] = ''
$data[180,0] = 'ACORN data is not of the right format. Only data generated with v2.1 (or later versions) is compatible.'
$data[180,1] = 'TBT'
$data[180,2] = 'deleted'
$data[181,0] = 'ACORN data is not of the right format. Only data generated with v2.1 is compatible.'
$data[181,1] = 'TBT'
$data[181,2] = 'deleted'
$data[182,0] = 'HAI point prevalence by type of ward'
$data[182,1] = 'ຄວາມຊຸກຊຸມ HAI ແບ່ງຕາມປະເພດຂອງພະແນກ'
$data[182,2] = 'deleted'
$data[183,0] = 'Select lab data format:'
$data[183,1] = 'ເລືອກຮູບແບບຂໍ້ມູມແລັບ'
$data[183,2] = 'deleted'

$ws.Range("A2:C185").Value = $data

Write-Output "Rewrote rows 2:185 (header row 1 untouched)."
